$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column values are stored as text (they are formatted strings like
# "25.006.98" or "1.000" that Excel would otherwise coerce to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.993.23"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.709.83"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "317.69"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "0.4046"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("D8").Value = "0.4087"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "1.482"
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").Value = "53.74"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "0.9989"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "0.08848"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "26.43"
$ws.Range("E13").Value = "  +6.75%  "
$ws.Range("D14").Value = "7.531"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "8.143"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "0.00001361"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "1.720.65"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").Value = "97.17"
$ws.Range("E18").Value = "  -3.30%  "
$ws.Range("D19").Value = "0.07160"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "21.24"
$ws.Range("E20").Value = "  +4.26%  "
$ws.Range("D21").Value = "7.291"
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "14.43"
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("D24").Value = "24.974.95"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "2.928"
$ws.Range("E25").Value = "  -7.00%  "
$ws.Range("D26").Value = "2.319"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "23.32"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("D28").Value = "6.209"
$ws.Range("E28").Value = "  +18.54%  "
$ws.Range("D29").Value = "167.14"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Value = "146.54"
$ws.Range("E30").Value = "  +4.75%  "
$ws.Range("D31").Value = "8.432"
$ws.Range("E31").Value = "  -9.24%  "
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.925.64"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("B33").Value = "WEMIXTOKEN"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "2.242"
$ws.Range("E33").Value = "  +14.03%  "
$ws.Range("D34").Value = "0.08900"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").Value = "0.03223"
$ws.Range("E35").Value = "  +6.52%  "
$ws.Range("D36").Value = "7.284"
$ws.Range("E36").Value = "  -7.63%  "
$ws.Range("D37").Value = "1.031"
$ws.Range("E37").Value = "  -6.13%  "
$ws.Range("D38").Value = "0.2866"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").Value = "0.8502"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").Value = "10.91"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("D41").Value = "0.09368"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "14.23"
$ws.Range("E42").Value = "  -3.34%  "
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("D44").Value = "17.44"
$ws.Range("E44").Value = "  +3.69%  "
$ws.Range("D45").Value = "2.719"
$ws.Range("E45").Value = "  +2.55%  "
$ws.Range("D46").Value = "0.7455"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "4.251"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "1.403"
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("D49").Value = "0.9998"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "142.28"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").Value = "0.08370"
$ws.Range("E51").Value = "  +3.13%  "
